$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the FAPs -> Lama1/Itga1 block (rows 2-6) with refreshed TPM-derived
# statistics, and append the new MuSCs -> Lama1/Itga1 block (rows 7-11).

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "FAPs"
$row2[0,1] = "Lama1"
$row2[0,2] = "Itga1"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 0.5587383333333333
$row2[0,7] = 1.676215
$row2[0,8] = 0.8486764927018626
$row2[0,9] = 0.8937587278261895
$row2[0,10] = 3
$row2[0,11] = 1
$row2[0,12] = 32.21373866666666
$row2[0,13] = 96.641216
$row2[0,14] = 0.6812298485843321
$row2[0,15] = 0.7117693664123
$row2[0,16] = 17.99905065304889
$row2[0,17] = 161.99145587744
$row2[0,18] = 0.5781437586203719
$row2[0,19] = 0.6361500834303102
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "FAPs"
$row3[0,1] = "Lama1"
$row3[0,2] = "Itga1"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 0.5587383333333333
$row3[0,7] = 1.676215
$row3[0,8] = 0.8486764927018626
$row3[0,9] = 0.8937587278261895
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 8.371752
$row3[0,13] = 25.115256
$row3[0,14] = 0.1770389772624213
$row3[0,15] = 0.184975630381169
$row3[0,16] = 4.677618759560001
$row3[0,17] = 42.09856883604
$row3[0,18] = 0.1502488182945965
$row3[0,19] = 0.1653235840883211
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "FAPs"
$row4[0,1] = "Lama1"
$row4[0,2] = "Itga1"
$row4[0,3] = "Inflammatory-Mac"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 0.5587383333333333
$row4[0,7] = 1.676215
$row4[0,8] = 0.8486764927018626
$row4[0,9] = 0.8937587278261895
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 0.5484013333333334
$row4[0,13] = 1.645204
$row4[0,14] = 0.01159714372603029
$row4[0,15] = 0.01211704340205096
$row4[0,16] = 0.3064128469844445
$row4[0,17] = 2.75771562286
$row4[0,18] = 0.009842223262766795
$row4[0,19] = 0.01082971329603179
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "FAPs"
$row5[0,1] = "Lama1"
$row5[0,2] = "Itga1"
$row5[0,3] = "MuSCs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 0.5587383333333333
$row5[0,7] = 1.676215
$row5[0,8] = 0.8486764927018626
$row5[0,9] = 0.8937587278261895
$row5[0,10] = 2
$row5[0,11] = 1
$row5[0,12] = 6.086836
$row5[0,13] = 12.173672
$row5[0,14] = 0.1287194389184112
$row5[0,15] = 0.08965995219214913
$row5[0,16] = 3.400948601913333
$row5[0,17] = 20.40569161148
$row5[0,18] = 0.1092411619638289
$row5[0,19] = 0.08013436480821218
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Lama1"
$row6[0,2] = "Itga1"
$row6[0,3] = "Resolving-Mac"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 0.5587383333333333
$row6[0,7] = 1.676215
$row6[0,8] = 0.8486764927018626
$row6[0,9] = 0.8937587278261895
$row6[0,10] = 2
$row6[0,11] = 0.6666666666666666
$row6[0,12] = 0.06689266666666667
$row6[0,13] = 0.200678
$row6[0,14] = 0.001414591508805173
$row6[0,15] = 0.001478007612330618
$row6[0,16] = 0.03737549708555556
$row6[0,17] = 0.33637947377
$row6[0,18] = 0.00120053056029861
$row6[0,19] = 0.001320982203314037
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "MuSCs"
$row7[0,1] = "Lama1"
$row7[0,2] = "Itga1"
$row7[0,3] = "ECs"
$row7[0,4] = 2
$row7[0,5] = 1
$row7[0,6] = 0.09962599999999999
$row7[0,7] = 0.199252
$row7[0,8] = 0.1513235072981373
$row7[0,9] = 0.1062412721738106
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 32.21373866666666
$row7[0,13] = 96.641216
$row7[0,14] = 0.6812298485843321
$row7[0,15] = 0.7117693664123
$row7[0,16] = 3.209325928405333
$row7[0,17] = 19.255955570432
$row7[0,18] = 0.1030860899639602
$row7[0,19] = 0.07561928298198987
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "MuSCs"
$row8[0,1] = "Lama1"
$row8[0,2] = "Itga1"
$row8[0,3] = "FAPs"
$row8[0,4] = 2
$row8[0,5] = 1
$row8[0,6] = 0.09962599999999999
$row8[0,7] = 0.199252
$row8[0,8] = 0.1513235072981373
$row8[0,9] = 0.1062412721738106
$row8[0,10] = 3
$row8[0,11] = 1
$row8[0,12] = 8.371752
$row8[0,13] = 25.115256
$row8[0,14] = 0.1770389772624213
$row8[0,15] = 0.184975630381169
$row8[0,16] = 0.834044164752
$row8[0,17] = 5.004264988512
$row8[0,18] = 0.02679015896782478
$row8[0,19] = 0.01965204629284796
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "MuSCs"
$row9[0,1] = "Lama1"
$row9[0,2] = "Itga1"
$row9[0,3] = "Inflammatory-Mac"
$row9[0,4] = 2
$row9[0,5] = 1
$row9[0,6] = 0.09962599999999999
$row9[0,7] = 0.199252
$row9[0,8] = 0.1513235072981373
$row9[0,9] = 0.1062412721738106
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 0.5484013333333334
$row9[0,13] = 1.645204
$row9[0,14] = 0.01159714372603029
$row9[0,15] = 0.01211704340205096
$row9[0,16] = 0.05463503123466667
$row9[0,17] = 0.327810187408
$row9[0,18] = 0.001754920463263492
$row9[0,19] = 0.001287330106019172
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "MuSCs"
$row10[0,1] = "Lama1"
$row10[0,2] = "Itga1"
$row10[0,3] = "MuSCs"
$row10[0,4] = 2
$row10[0,5] = 1
$row10[0,6] = 0.09962599999999999
$row10[0,7] = 0.199252
$row10[0,8] = 0.1513235072981373
$row10[0,9] = 0.1062412721738106
$row10[0,10] = 2
$row10[0,11] = 1
$row10[0,12] = 6.086836
$row10[0,13] = 12.173672
$row10[0,14] = 0.1287194389184112
$row10[0,15] = 0.08965995219214913
$row10[0,16] = 0.606407123336
$row10[0,17] = 2.425628493344
$row10[0,18] = 0.01947827695458235
$row10[0,19] = 0.00952558738393696
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = "MuSCs"
$row11[0,1] = "Lama1"
$row11[0,2] = "Itga1"
$row11[0,3] = "Resolving-Mac"
$row11[0,4] = 2
$row11[0,5] = 1
$row11[0,6] = 0.09962599999999999
$row11[0,7] = 0.199252
$row11[0,8] = 0.1513235072981373
$row11[0,9] = 0.1062412721738106
$row11[0,10] = 2
$row11[0,11] = 0.6666666666666666
$row11[0,12] = 0.06689266666666667
$row11[0,13] = 0.200678
$row11[0,14] = 0.001414591508805173
$row11[0,15] = 0.001478007612330618
$row11[0,16] = 0.006664248809333333
$row11[0,17] = 0.039985492856
$row11[0,18] = 0.0002140609485065627
$row11[0,19] = 0.0001570254090165811
$ws.Range("A11:T11").Value = $row11
